$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change landmark_placement value from "manual" to "manually"
$ws.Range("J2").Value = "manually"

# Add new column M: "verbose" header (bold, bordered like other headers) with value TRUE
$ws.Range("M1").Value = "verbose"
$ws.Range("M1").Style = $ws.Range("L1").Style

$ws.Range("M2").Value = $true
$ws.Range("M2").Style = $ws.Range("L2").Style

# Update selection to match the target state
$ws.Range("E5").Select()
